$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Razon social" (column E) text that got its commas mangled into periods
#     by the same (buggy) float-normalization pass the scraper used on numbers. ---
$ws.Range("E41").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E103").Value = 'FERNANDEZ. MARIO HUGO'
$ws.Range("E105").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E173").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E178").Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range("E196").Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'

# --- Fix "Importe" (column H): scraped floats were formatted with a thousands-dot
#     and a decimal-comma ("2.720,00"); normalize to a plain decimal-dot string
#     ("2720.00") while keeping the cells as text. ---
$ws.Range("H2:H241").NumberFormat = "@"
$ws.Range("H2").Value = '2720.00'
$ws.Range("H3").Value = '54000.00'
$ws.Range("H4").Value = '41350.00'
$ws.Range("H5").Value = '51550.00'
$ws.Range("H6").Value = '1301781.00'
$ws.Range("H7").Value = '10500.00'
$ws.Range("H8").Value = '734.98'
$ws.Range("H9").Value = '798.00'
$ws.Range("H10").Value = '745.00'
$ws.Range("H11").Value = '3025.00'
$ws.Range("H12").Value = '85000.00'
$ws.Range("H13").Value = '14340.00'
$ws.Range("H14").Value = '5483.98'
$ws.Range("H15").Value = '4098.38'
$ws.Range("H16").Value = '106675.00'
$ws.Range("H17").Value = '629237.84'
$ws.Range("H18").Value = '10659.40'
$ws.Range("H19").Value = '8390.00'
$ws.Range("H20").Value = '7990.80'
$ws.Range("H21").Value = '9455.00'
$ws.Range("H22").Value = '14131.60'
$ws.Range("H23").Value = '8865.00'
$ws.Range("H24").Value = '2666.00'
$ws.Range("H25").Value = '6021.00'
$ws.Range("H26").Value = '12239.32'
$ws.Range("H27").Value = '532.00'
$ws.Range("H28").Value = '3395.03'
$ws.Range("H29").Value = '12049.99'
$ws.Range("H30").Value = '6376.51'
$ws.Range("H31").Value = '3700.00'
$ws.Range("H32").Value = '1250.00'
$ws.Range("H33").Value = '50064.97'
$ws.Range("H34").Value = '2400.00'
$ws.Range("H35").Value = '148.93'
$ws.Range("H36").Value = '926.90'
$ws.Range("H37").Value = '31849.40'
$ws.Range("H38").Value = '18321.18'
$ws.Range("H39").Value = '32806.65'
$ws.Range("H40").Value = '3995.21'
$ws.Range("H41").Value = '700.00'
$ws.Range("H42").Value = '4150.35'
$ws.Range("H43").Value = '31.38'
$ws.Range("H44").Value = '64.74'
$ws.Range("H45").Value = '13438.00'
$ws.Range("H46").Value = '315.00'
$ws.Range("H47").Value = '496.50'
$ws.Range("H48").Value = '234.70'
$ws.Range("H49").Value = '317.64'
$ws.Range("H50").Value = '1152682.00'
$ws.Range("H51").Value = '24416.20'
$ws.Range("H52").Value = '833.28'
$ws.Range("H53").Value = '270.00'
$ws.Range("H54").Value = '12039.00'
$ws.Range("H55").Value = '22797.00'
$ws.Range("H56").Value = '17941.77'
$ws.Range("H57").Value = '9354.80'
$ws.Range("H58").Value = '402.14'
$ws.Range("H59").Value = '70.00'
$ws.Range("H60").Value = '923.65'
$ws.Range("H61").Value = '4407.61'
$ws.Range("H62").Value = '1219.94'
$ws.Range("H63").Value = '4867.80'
$ws.Range("H64").Value = '13500.00'
$ws.Range("H65").Value = '11010.00'
$ws.Range("H66").Value = '3502.80'
$ws.Range("H67").Value = '4180.08'
$ws.Range("H68").Value = '180.84'
$ws.Range("H69").Value = '3.40'
$ws.Range("H70").Value = '99800.76'
$ws.Range("H71").Value = '1078.10'
$ws.Range("H72").Value = '5286.77'
$ws.Range("H73").Value = '1539.29'
$ws.Range("H74").Value = '443.54'
$ws.Range("H75").Value = '976.00'
$ws.Range("H76").Value = '121.00'
$ws.Range("H77").Value = '1940.20'
$ws.Range("H78").Value = '9926.65'
$ws.Range("H79").Value = '38267.45'
$ws.Range("H80").Value = '3500.00'
$ws.Range("H81").Value = '6300.00'
$ws.Range("H82").Value = '492.71'
$ws.Range("H83").Value = '1165.00'
$ws.Range("H84").Value = '989.00'
$ws.Range("H85").Value = '1924.30'
$ws.Range("H86").Value = '12608.40'
$ws.Range("H87").Value = '242.20'
$ws.Range("H88").Value = '278.86'
$ws.Range("H89").Value = '2076.51'
$ws.Range("H90").Value = '4070.00'
$ws.Range("H91").Value = '2440.00'
$ws.Range("H92").Value = '165.00'
$ws.Range("H93").Value = '440.00'
$ws.Range("H94").Value = '500.00'
$ws.Range("H95").Value = '140.00'
$ws.Range("H96").Value = '400.00'
$ws.Range("H97").Value = '21252.00'
$ws.Range("H98").Value = '32051.00'
$ws.Range("H99").Value = '9700.00'
$ws.Range("H100").Value = '498.00'
$ws.Range("H101").Value = '35477.00'
$ws.Range("H102").Value = '1176.00'
$ws.Range("H103").Value = '1678.00'
$ws.Range("H104").Value = '16680.00'
$ws.Range("H105").Value = '2422.00'
$ws.Range("H106").Value = '468.67'
$ws.Range("H107").Value = '999.96'
$ws.Range("H108").Value = '52.85'
$ws.Range("H109").Value = '1300.00'
$ws.Range("H110").Value = '7294.78'
$ws.Range("H111").Value = '120.00'
$ws.Range("H112").Value = '12000.00'
$ws.Range("H113").Value = '345761.96'
$ws.Range("H114").Value = '53261.49'
$ws.Range("H115").Value = '1360.00'
$ws.Range("H116").Value = '240.26'
$ws.Range("H117").Value = '9163.91'
$ws.Range("H118").Value = '1152.50'
$ws.Range("H119").Value = '810.00'
$ws.Range("H120").Value = '490.00'
$ws.Range("H121").Value = '2340.00'
$ws.Range("H122").Value = '17152.56'
$ws.Range("H123").Value = '5140.70'
$ws.Range("H124").Value = '50.00'
$ws.Range("H125").Value = '356.00'
$ws.Range("H126").Value = '250.00'
$ws.Range("H127").Value = '9.12'
$ws.Range("H128").Value = '95.00'
$ws.Range("H129").Value = '4500.00'
$ws.Range("H130").Value = '380.00'
$ws.Range("H131").Value = '4098.00'
$ws.Range("H132").Value = '2125.68'
$ws.Range("H133").Value = '568.00'
$ws.Range("H134").Value = '494.50'
$ws.Range("H135").Value = '60.00'
$ws.Range("H136").Value = '113.44'
$ws.Range("H137").Value = '5732.84'
$ws.Range("H138").Value = '5407.00'
$ws.Range("H139").Value = '960.00'
$ws.Range("H140").Value = '682.00'
$ws.Range("H141").Value = '380.00'
$ws.Range("H142").Value = '36000.00'
$ws.Range("H143").Value = '2450.00'
$ws.Range("H144").Value = '4645.00'
$ws.Range("H145").Value = '21960.00'
$ws.Range("H146").Value = '83564.00'
$ws.Range("H147").Value = '4328.83'
$ws.Range("H148").Value = '18.56'
$ws.Range("H149").Value = '4150.00'
$ws.Range("H150").Value = '1490.00'
$ws.Range("H151").Value = '5521.50'
$ws.Range("H152").Value = '1440.00'
$ws.Range("H153").Value = '3740.00'
$ws.Range("H154").Value = '800.00'
$ws.Range("H155").Value = '2500.00'
$ws.Range("H156").Value = '96800.00'
$ws.Range("H157").Value = '1500.00'
$ws.Range("H158").Value = '1566.08'
$ws.Range("H159").Value = '1360.50'
$ws.Range("H160").Value = '3877.50'
$ws.Range("H161").Value = '110.01'
$ws.Range("H162").Value = '962.40'
$ws.Range("H163").Value = '2589.99'
$ws.Range("H164").Value = '292.10'
$ws.Range("H165").Value = '700000.00'
$ws.Range("H166").Value = '30580.80'
$ws.Range("H167").Value = '850.00'
$ws.Range("H168").Value = '12298.65'
$ws.Range("H169").Value = '500.00'
$ws.Range("H170").Value = '32390.40'
$ws.Range("H171").Value = '488.00'
$ws.Range("H172").Value = '300.00'
$ws.Range("H173").Value = '2275.00'
$ws.Range("H174").Value = '700.00'
$ws.Range("H175").Value = '492.20'
$ws.Range("H176").Value = '4680.00'
$ws.Range("H177").Value = '64912.00'
$ws.Range("H178").Value = '380.00'
$ws.Range("H179").Value = '42764.00'
$ws.Range("H180").Value = '1634.57'
$ws.Range("H181").Value = '956.00'
$ws.Range("H182").Value = '4340.00'
$ws.Range("H183").Value = '6340.00'
$ws.Range("H184").Value = '165.92'
$ws.Range("H185").Value = '187.50'
$ws.Range("H186").Value = '220.20'
$ws.Range("H187").Value = '1360.00'
$ws.Range("H188").Value = '506.00'
$ws.Range("H189").Value = '81.81'
$ws.Range("H190").Value = '336.34'
$ws.Range("H191").Value = '9823.88'
$ws.Range("H192").Value = '3486.00'
$ws.Range("H193").Value = '80709.00'
$ws.Range("H194").Value = '4559.36'
$ws.Range("H195").Value = '165.00'
$ws.Range("H196").Value = '2175.00'
$ws.Range("H197").Value = '1558.28'
$ws.Range("H198").Value = '1955.67'
$ws.Range("H199").Value = '3347.52'
$ws.Range("H200").Value = '490.00'
$ws.Range("H201").Value = '4095.00'
$ws.Range("H202").Value = '1908.40'
$ws.Range("H203").Value = '70.00'
$ws.Range("H204").Value = '2961.00'
$ws.Range("H205").Value = '745.00'
$ws.Range("H206").Value = '1420.00'
$ws.Range("H207").Value = '50217.30'
$ws.Range("H208").Value = '982.00'
$ws.Range("H209").Value = '359.07'
$ws.Range("H210").Value = '4077.25'
$ws.Range("H211").Value = '9450.00'
$ws.Range("H212").Value = '2786.90'
$ws.Range("H213").Value = '795867.82'
$ws.Range("H214").Value = '1524.60'
$ws.Range("H215").Value = '708.26'
$ws.Range("H216").Value = '1600.00'
$ws.Range("H217").Value = '51100.00'
$ws.Range("H218").Value = '2120.00'
$ws.Range("H219").Value = '4980.00'
$ws.Range("H220").Value = '34774.00'
$ws.Range("H221").Value = '5300.00'
$ws.Range("H222").Value = '179336.00'
$ws.Range("H223").Value = '189240.00'
$ws.Range("H224").Value = '42800.00'
$ws.Range("H225").Value = '108020.00'
$ws.Range("H226").Value = '107650.00'
$ws.Range("H227").Value = '35000.00'
$ws.Range("H228").Value = '59500.00'
$ws.Range("H229").Value = '244344.00'
$ws.Range("H230").Value = '137948.00'
$ws.Range("H231").Value = '268304.00'
$ws.Range("H232").Value = '191553.75'
$ws.Range("H233").Value = '7900.00'
$ws.Range("H234").Value = '360580.00'
$ws.Range("H235").Value = '25900.00'
$ws.Range("H236").Value = '4000.00'
$ws.Range("H237").Value = '14215.50'
$ws.Range("H238").Value = '511500.00'
$ws.Range("H239").Value = '854.00'
$ws.Range("H240").Value = '2800.00'
$ws.Range("H241").Value = '1550.00'
